$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-5 from 45208 to 45212
$ws.Range("C2:C5").Value = 45212
